# "#1 Added labels to graphic"
#
# texts.xlsx has two sheets:
#   - "Typography" : table of font/typography definitions (B3:I100)
#   - "Translation" : table of text-id -> translation strings (B3:I799)
#
# This change adds a new Typography entry ("Typography_01", a smaller size
# variant of the existing "angsai.ttf" typography used for the graphic's
# axis-value labels), and reworks a couple of the leftover
# "SingleUseIdN" placeholder rows on the Translation sheet into two real,
# named text ids ("TextId2" / "TextId3") that carry the new "<currXValue>"
# / "<currYValue>" graphic-label placeholders, clearing out the remaining
# now-unused placeholder rows.

$wb = $excel.ActiveWorkbook
$typography = $wb.Worksheets.Item("Typography")
$translation = $wb.Worksheets.Item("Translation")

# ---------------------------------------------------------------------
# Typography sheet: shrink the "Typography_00" font size and give it a
# widget wildcard-character set, then add a new "Typography_01" row
# directly below it (same font, smaller size) for the graphic labels.
# ---------------------------------------------------------------------

$typography.Range("D5").Value = 15
$typography.Range("J5").Value = "-., 0123456789"

$typography.Range("B6").Value = "Typography_01"
$typography.Range("C6").Value = "angsai.ttf"
$typography.Range("D6").Value = 20
$typography.Range("E6").Value = 4
$typography.Range("F6").Value = "?"
$typography.Range("G6").Value = ""
$typography.Range("H6").Value = "0-9"
$typography.Range("J6").Value = ""

# B6:E6 would otherwise inherit the column's (non-default) style from the
# <cols> definition; reset them (and the two blank marker cells) back to
# the un-styled "Normal" look the rest of the table rows use.
$typography.Range("B6:E6").Style = "Normal"
$typography.Range("G6").Style = "Normal"
$typography.Range("J6").Style = "Normal"

# ---------------------------------------------------------------------
# Translation sheet: point the "SingleUseId2" row at the new
# "Typography_01" typography, then turn the old SingleUseId5/6 rows into
# the new "TextId2" / "TextId3" graphic-value-label entries, and clear
# out the remaining leftover SingleUseId7-10 rows.
# ---------------------------------------------------------------------

$translation.Range("C6").Value = "Typography_01"

$translation.Range("B9").Value = "TextId2"
$translation.Range("C9").Value = "Typography_00"
$translation.Range("D9").Value = "Center"
$translation.Range("E9").Value = "<number>"

$translation.Range("B10").Value = "TextId3"
$translation.Range("C10").Value = "Typography_01"
$translation.Range("D10").Value = "Center"
$translation.Range("E10").Value = "<number>"

$translation.Range("B11:F11").ClearContents()
$translation.Range("B12:F12").ClearContents()
$translation.Range("B13:F13").ClearContents()
$translation.Range("B14:F14").ClearContents()
